$wb = $excel.ActiveWorkbook

$wsCases = $wb.Worksheets.Item("TestCases")
$wsTest  = $wb.Worksheets.Item("Test")

# --- Update the existing "Scenario Radio Button" (TC_30) result rows -------
# These cells already existed but held placeholder text; fill in the real
# recorded results now that the Radio Button page has been tested.
$wsTest.Range("B54").Value = "Header of RadioButton"
$wsTest.Range("B55").Value = "Radio Button"
$wsTest.Range("C55").Value = "Radio Button"

# --- TC_31: elements on the Radio Button page -------------------------------
$wsTest.Range("A48:D48").Copy($wsTest.Range("A56"))
$wsTest.Range("A56").Value = "TC_31"
$wsTest.Range("B56").Value = "Elements Of RadioButtonPage"

$wsTest.Range("B57").Value = "Do you like the site?"
$wsTest.Range("C57").Value = "Do you like the site?"

$wsTest.Range("B58").Value = "Yes"
$wsTest.Range("C58").Value = "Yes"

$wsTest.Range("B59").Value = "Impressive"
$wsTest.Range("C59").Value = "Impressive"

$wsTest.Range("B60").Value = "No"
$wsTest.Range("C60").Value = "No"

# --- TC_32: selecting the "Yes" radio button --------------------------------
$wsTest.Range("A48:D48").Copy($wsTest.Range("A61"))
$wsTest.Range("A61").Value = "TC_32"
$wsTest.Range("B61").Value = "Text Selected Yes Radio Button "

$wsTest.Range("B62").Value = "You have selected Yes"
$wsTest.Range("C62").Value = "You have selected Yes"

# --- TC_33: selecting the "Impressive" radio button -------------------------
$wsTest.Range("A48:D48").Copy($wsTest.Range("A63"))
$wsTest.Range("A63").Value = "TC_33"
$wsTest.Range("B63").Value = "Text Selected Yes Radio Button "

$wsTest.Range("B64").Value = "You have selected Impressive"
$wsTest.Range("C64").Value = "You have selected Impressive"

# --- View state: last selection on TestCases was C241 ----------------------
$null = $wsCases.Range("C241").Select()

# --- Make "Test" the active sheet, with its last selection at B76 ----------
$wsTest.Activate()
$null = $wsTest.Range("B76").Select()
